# Auto-generated edit script: Add data for 2025-09-08
# Updates year-2025 (column L) violent crime counts across Citywide, By Neighborhood,
# and individual neighborhood sheets in the workbook.

$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, [string]$cellRef, $newValue)
    $ws.Range($cellRef).Value = $newValue
}

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
Set-CellValue $ws "L2" 4636
Set-CellValue $ws "L3" 4994
Set-CellValue $ws "K4" 1780
Set-CellValue $ws "L4" 1231
Set-CellValue $ws "L5" 290
Set-CellValue $ws "L6" 4254
Set-CellValue $ws "K7" 27572
Set-CellValue $ws "L7" 15405

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
Set-CellValue $ws "L2" 128
Set-CellValue $ws "L7" 506
Set-CellValue $ws "L8" 1023
Set-CellValue $ws "L11" 247
Set-CellValue $ws "L14" 78
Set-CellValue $ws "L19" 424
Set-CellValue $ws "L20" 390
Set-CellValue $ws "L23" 166
Set-CellValue $ws "L27" 136
Set-CellValue $ws "L29" 853
Set-CellValue $ws "L31" 150
Set-CellValue $ws "L32" 20
Set-CellValue $ws "L33" 701
Set-CellValue $ws "L34" 91
Set-CellValue $ws "L37" 566
Set-CellValue $ws "L42" 502
Set-CellValue $ws "L45" 28
Set-CellValue $ws "L47" 108
Set-CellValue $ws "L48" 200
Set-CellValue $ws "L52" 311
Set-CellValue $ws "L54" 322
Set-CellValue $ws "L55" 146
Set-CellValue $ws "L57" 55
Set-CellValue $ws "L60" 97
Set-CellValue $ws "L61" 17
Set-CellValue $ws "K63" 168
Set-CellValue $ws "L63" 45
Set-CellValue $ws "L65" 299
Set-CellValue $ws "L67" 530
Set-CellValue $ws "L73" 122
Set-CellValue $ws "L76" 240
Set-CellValue $ws "L78" 203
Set-CellValue $ws "L79" 408
Set-CellValue $ws "L83" 338
Set-CellValue $ws "L84" 149
Set-CellValue $ws "L85" 792
Set-CellValue $ws "L86" 114
Set-CellValue $ws "L87" 45
Set-CellValue $ws "L88" 165
Set-CellValue $ws "L89" 219
Set-CellValue $ws "L91" 208
Set-CellValue $ws "L93" 80
Set-CellValue $ws "L94" 193
Set-CellValue $ws "L95" 209
Set-CellValue $ws "L97" 132
Set-CellValue $ws "L99" 264
Set-CellValue $ws "L100" 24
Set-CellValue $ws "K101" 27572
Set-CellValue $ws "L101" 15405

# Sheet 3: Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
Set-CellValue $ws "L6" 19
Set-CellValue $ws "L7" 78

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
Set-CellValue $ws "L2" 171
Set-CellValue $ws "L7" 506

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
Set-CellValue $ws "L3" 77
Set-CellValue $ws "L7" 247

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item("Uptown")
Set-CellValue $ws "L6" 60
Set-CellValue $ws "L7" 219

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item("South Shore")
Set-CellValue $ws "L2" 239
Set-CellValue $ws "L3" 320
Set-CellValue $ws "L6" 167
Set-CellValue $ws "L7" 792

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item("Little Village")
Set-CellValue $ws "L3" 98
Set-CellValue $ws "L7" 311

# Sheet 12: Austin
$ws = $wb.Worksheets.Item("Austin")
Set-CellValue $ws "L2" 294
Set-CellValue $ws "L3" 345
Set-CellValue $ws "L5" 37
Set-CellValue $ws "L7" 1023

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
Set-CellValue $ws "L3" 133
Set-CellValue $ws "L6" 79
Set-CellValue $ws "L7" 338

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
Set-CellValue $ws "L2" 191
Set-CellValue $ws "L3" 237
Set-CellValue $ws "L4" 42
Set-CellValue $ws "L6" 216
Set-CellValue $ws "L7" 701

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
Set-CellValue $ws "L3" 67
Set-CellValue $ws "L7" 209

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
Set-CellValue $ws "L6" 154
Set-CellValue $ws "L7" 566

# Sheet 17: New City
$ws = $wb.Worksheets.Item("New City")
Set-CellValue $ws "L2" 106
Set-CellValue $ws "L3" 95
Set-CellValue $ws "L7" 299

# Sheet 18: Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
Set-CellValue $ws "L2" 72
Set-CellValue $ws "L6" 58
Set-CellValue $ws "L7" 264

# Sheet 20: Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
Set-CellValue $ws "L3" 38
Set-CellValue $ws "L7" 150

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
Set-CellValue $ws "L2" 154
Set-CellValue $ws "L6" 121
Set-CellValue $ws "L7" 530

# Sheet 22: South Deering
$ws = $wb.Worksheets.Item("South Deering")
Set-CellValue $ws "L6" 42
Set-CellValue $ws "L7" 149

# Sheet 24: Loop
$ws = $wb.Worksheets.Item("Loop")
Set-CellValue $ws "L3" 78
Set-CellValue $ws "L7" 322

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item("Englewood")
Set-CellValue $ws "L3" 321
Set-CellValue $ws "L4" 41
Set-CellValue $ws "L6" 221
Set-CellValue $ws "L7" 853

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item("Lake View")
Set-CellValue $ws "L3" 49
Set-CellValue $ws "L6" 86
Set-CellValue $ws "L7" 200

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item("Chatham")
Set-CellValue $ws "L2" 150
Set-CellValue $ws "L3" 131
Set-CellValue $ws "L6" 122
Set-CellValue $ws "L7" 424

# Sheet 29: River North
$ws = $wb.Worksheets.Item("River North")
Set-CellValue $ws "L3" 44
Set-CellValue $ws "L7" 240

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
Set-CellValue $ws "L2" 144
Set-CellValue $ws "L3" 169
Set-CellValue $ws "L4" 36
Set-CellValue $ws "L6" 140
Set-CellValue $ws "L7" 502

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
Set-CellValue $ws "L2" 55
Set-CellValue $ws "L3" 66
Set-CellValue $ws "L7" 203

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
Set-CellValue $ws "L2" 46
Set-CellValue $ws "L3" 48
Set-CellValue $ws "L6" 43
Set-CellValue $ws "L7" 146

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item("Douglas")
Set-CellValue $ws "L3" 61
Set-CellValue $ws "L7" 166

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
Set-CellValue $ws "L3" 91
Set-CellValue $ws "L7" 208

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item("Roseland")
Set-CellValue $ws "L2" 135
Set-CellValue $ws "L6" 86
Set-CellValue $ws "L7" 408

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
Set-CellValue $ws "L3" 129
Set-CellValue $ws "L7" 390

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
Set-CellValue $ws "L6" 25
Set-CellValue $ws "L7" 80

# Sheet 49: Wrigleyville
$ws = $wb.Worksheets.Item("Wrigleyville")
Set-CellValue $ws "L4" 2
Set-CellValue $ws "L6" 13
Set-CellValue $ws "L7" 24

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
Set-CellValue $ws "L6" 30
Set-CellValue $ws "L7" 91

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item("West Loop")
Set-CellValue $ws "L3" 43
Set-CellValue $ws "L7" 193

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
Set-CellValue $ws "L2" 41
Set-CellValue $ws "L7" 108

# Sheet 62: Portage Park
$ws = $wb.Worksheets.Item("Portage Park")
Set-CellValue $ws "L2" 44
Set-CellValue $ws "L7" 122

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item("Albany Park")
Set-CellValue $ws "L4" 10
Set-CellValue $ws "L7" 128

# Sheet 65: West Town
$ws = $wb.Worksheets.Item("West Town")
Set-CellValue $ws "L2" 29
Set-CellValue $ws "L7" 132

# Sheet 68: United Center
$ws = $wb.Worksheets.Item("United Center")
Set-CellValue $ws "L5" 6
Set-CellValue $ws "L7" 165

# Sheet 69: Galewood
$ws = $wb.Worksheets.Item("Galewood")
Set-CellValue $ws "L2" 14
Set-CellValue $ws "L7" 20

# Sheet 71: Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
Set-CellValue $ws "L2" 36
Set-CellValue $ws "L7" 136

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item("Streeterville")
Set-CellValue $ws "L4" 62
Set-CellValue $ws "L7" 114

# Sheet 77: Mckinley Park
$ws = $wb.Worksheets.Item("Mckinley Park")
Set-CellValue $ws "L5" 1
Set-CellValue $ws "L7" 55

# Sheet 78: Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
Set-CellValue $ws "L3" 35
Set-CellValue $ws "L7" 97

# Sheet 85: Jackson Park
$ws = $wb.Worksheets.Item("Jackson Park")
Set-CellValue $ws "L6" 6
Set-CellValue $ws "L7" 28

# Sheet 92: Ukrainian Village
$ws = $wb.Worksheets.Item("Ukrainian Village")
Set-CellValue $ws "L2" 13
Set-CellValue $ws "L7" 45

# Sheet 93: Mount Greenwood
$ws = $wb.Worksheets.Item("Mount Greenwood")
Set-CellValue $ws "L6" 5
Set-CellValue $ws "L7" 17
